$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell for new column F, matching style of existing header cells (A1:E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

$timestamps = @(
    "2021-10-05 10:51:11.460882",
    "2021-10-05 10:51:11.460895",
    "2021-10-05 10:51:11.460899",
    "2021-10-05 10:51:11.460902",
    "2021-10-05 10:51:11.460906",
    "2021-10-05 10:51:11.460909",
    "2021-10-05 10:51:11.460912",
    "2021-10-05 10:51:11.460915",
    "2021-10-05 10:51:11.460918",
    "2021-10-05 10:51:11.460921",
    "2021-10-05 10:51:11.460924",
    "2021-10-05 10:51:11.460927",
    "2021-10-05 10:51:11.460930",
    "2021-10-05 10:51:11.460933",
    "2021-10-05 10:51:11.460936",
    "2021-10-05 10:51:11.460940",
    "2021-10-05 10:51:11.460943",
    "2021-10-05 10:51:11.460946",
    "2021-10-05 10:51:11.460949",
    "2021-10-05 10:51:11.460952",
    "2021-10-05 10:51:11.460955",
    "2021-10-05 10:51:11.460958",
    "2021-10-05 10:51:11.460961",
    "2021-10-05 10:51:11.460964",
    "2021-10-05 10:51:11.460968",
    "2021-10-05 10:51:11.460971"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
